$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'42.964.47"
$ws.Range('E2').Value = "  +0.00%  "
$ws.Range('D3').Value = "'2.211.59"
$ws.Range('E3').Value = "  -1.00%  "
$ws.Range('E4').Value = "  -0.11%  "
$ws.Range('D5').Value = "'257.56"
$ws.Range('E5').Value = "  +2.45%  "
$ws.Range('D6').Value = "'0.614"
$ws.Range('E6').Value = "  +0.00%  "
$ws.Range('D7').Value = "'76.94"
$ws.Range('E7').Value = "  +2.81%  "
$ws.Range('E8').Value = "  -0.03%  "
$ws.Range('D9').Value = "'0.596"
$ws.Range('E9').Value = "  -0.11%  "
$ws.Range('D10').Value = "'42.99"
$ws.Range('E10').Value = "  +3.77%  "
$ws.Range('D11').Value = "'0.0909"
$ws.Range('E11').Value = "  -1.60%  "
$ws.Range('D12').Value = "'6.99"
$ws.Range('E12').Value = "  +1.51%  "
$ws.Range('D14').Value = "'2.542.44"
$ws.Range('E14').Value = "  -1.12%  "
$ws.Range('D15').Value = "'14.43"
$ws.Range('E15').Value = "  -0.36%  "
$ws.Range('D16').Value = "'2.209.84"
$ws.Range('E16').Value = "  -1.24%  "
$ws.Range('D17').Value = "'0.785"
$ws.Range('E17').Value = "  -0.21%  "
$ws.Range('D18').Value = "'42.940.53"
$ws.Range('E18').Value = "  +0.20%  "
$ws.Range('E19').Value = "  +0.48%  "
$ws.Range('D20').Value = "'71.19"
$ws.Range('E20').Value = "  -0.05%  "
$ws.Range('D21').Value = "'5.99"
$ws.Range('E21').Value = "  +0.90%  "
$ws.Range('D22').Value = "'2.36"
$ws.Range('E22').Value = "  +8.97%  "
$ws.Range('D23').Value = "'230.27"
$ws.Range('E23').Value = "  +0.21%  "
$ws.Range('D24').Value = "'9.24"
$ws.Range('E24').Value = "  -1.99%  "
$ws.Range('E25').Value = "  -0.03%  "
$ws.Range('D26').Value = "'42.50"
$ws.Range('E26').Value = "  +9.23%  "
$ws.Range('D27').Value = "'10.76"
$ws.Range('E27').Value = "  +0.79%  "
$ws.Range('D28').Value = "'3.34"
$ws.Range('E28').Value = "  -2.81%  "
$ws.Range('E29').Value = "  -0.57%  "
$ws.Range('E30').Value = "  +3.21%  "
$ws.Range('D31').Value = "'172.72"
$ws.Range('E31').Value = "  +0.82%  "
$ws.Range('D32').Value = "'20.36"
$ws.Range('E32').Value = "  +0.90%  "
$ws.Range('D33').Value = "'0.0868"
$ws.Range('E33').Value = "  +9.36%  "
$ws.Range('D34').Value = "'5.25"
$ws.Range('E34').Value = "  +0.78%  "
$ws.Range('E35').Value = "  +0.54%  "
$ws.Range('D36').Value = "'0.0363"
$ws.Range('E36').Value = "  +11.03%  "
$ws.Range('D37').Value = "'0.107"
$ws.Range('E37').Value = "  -3.72%  "
$ws.Range('D38').Value = "'4.39"
$ws.Range('E38').Value = "  -0.86%  "
$ws.Range('D39').Value = "'12.95"
$ws.Range('E39').Value = "  +2.92%  "
$ws.Range('D40').Value = "'2.92"
$ws.Range('E40').Value = "  +19.82%  "
$ws.Range('D41').Value = "'2.12"
$ws.Range('E41').Value = "  +1.13%  "
$ws.Range('D42').Value = "'0.203"
$ws.Range('E42').Value = "  -0.93%  "
$ws.Range('D43').Value = "'61.26"
$ws.Range('E43').Value = "  +3.07%  "
$ws.Range('D44').Value = "'5.28"
$ws.Range('E44').Value = "  -1.84%  "
$ws.Range('D45').Value = "'103.26"
$ws.Range('E45').Value = "  +0.55%  "
$ws.Range('D46').Value = "'8.50"
$ws.Range('E46').Value = "  -2.08%  "
$ws.Range('D47').Value = "'0.470"
$ws.Range('E47').Value = "  -2.06%  "
$ws.Range('E48').Value = "  -1.53%  "
$ws.Range('D49').Value = "'1.12"
$ws.Range('E49').Value = "  +0.18%  "
$ws.Range('E50').Value = "  -0.59%  "
$ws.Range('E51').Value = "  +21.55%  "
